$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.143.44'
$ws.Range("E2").Value = '  +3.09%  '

$ws.Range("D3").Value = '3.811.23'
$ws.Range("E3").Value = '  +1.05%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '706.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +11.68%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.50'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.69%  '

$ws.Range("D7").Value = '3.809.69'
$ws.Range("E7").Value = '  +1.07%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.526'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.36%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.164'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.60%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.43'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +9.94%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.463'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.48%  '

$ws.Range("E13").Value = '  +7.66%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.38'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.64%  '

$ws.Range("D15").Value = '4.449.56'
$ws.Range("E15").Value = '  +1.02%  '

$ws.Range("D16").Value = '3.812.57'
$ws.Range("E16").Value = '  +1.18%  '

$ws.Range("D17").Value = '71.163.91'
$ws.Range("E17").Value = '  +3.12%  '

$ws.Range("E18").Value = '  +2.14%  '

$ws.Range("E19").Value = '  +3.45%  '

$ws.Range("E20").Value = '  +0.44%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.23'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +18.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '484.17'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.74%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.719'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.34%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.33%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000147'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.77%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.54'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.72%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.75%  '

$ws.Range("E28").Value = '  +4.06%  '

$ws.Range("D29").Value = '3.960.12'
$ws.Range("E29").Value = '  +0.89%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.15%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.08'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +15.10%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.61'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.68%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.31'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.70%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.75'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.69%  '

$ws.Range("E35").Value = '  +2.13%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.29'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.40%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.02%  '

$ws.Range("D38").Value = '3.760.01'
$ws.Range("E38").Value = '  +1.01%  '

$ws.Range("E39").Value = '  +3.70%  '

$ws.Range("E40").Value = '  +7.07%  '

$ws.Range("E41").Value = '  +4.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.27'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +14.62%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.000340'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +29.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.975'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '45.83'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.47%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '49.45'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.41%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '160.48'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.62%  '

$ws.Range("E50").Value = '  -0.06%  '

$ws.Range("E51").Value = '  +2.95%  '
